# Updated cryptos list on Sat Sep  9 20:33:23 UTC 2023 with GitHub Actions
# Refreshes price/volume columns for the crypto tracker sheet.
# Numeric-looking text prices (e.g. "216.24") are written with a leading
# apostrophe so Excel keeps them as text, matching the sheet's existing
# inline-string formatting instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.038.96'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '1.644.17'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('D4').Value = '''1.01'
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').Value = '''216.24'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('D6').Value = '''0.508'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '''1.01'
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('E8').Value = '  +0.47%  '
$ws.Range('D9').Value = '''0.0639'
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('D10').Value = '''19.59'
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('D11').Value = '''0.0798'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '''4.28'
$ws.Range('E12').Value = '  +0.79%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.671.02'
$ws.Range('E13').Value = '  +1.76%  '
$ws.Range('D14').Value = '''0.545'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').Value = '0.0₃0764'
$ws.Range('E15').Value = '  +1.16%  '
$ws.Range('D16').Value = '''63.46'
$ws.Range('E16').Value = '  +1.39%  '
$ws.Range('D17').Value = '26.062.36'
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').Value = '''194.60'
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('D21').Value = '''9.94'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = '''6.21'
$ws.Range('E22').Value = '  -1.06%  '
$ws.Range('E23').Value = '  +4.59%  '
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('D26').Value = '''143.27'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('E27').Value = '  +0.77%  '
$ws.Range('E28').Value = '  +0.58%  '
$ws.Range('E29').Value = '  +0.65%  '
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('D31').Value = '''3.30'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').Value = '''3.27'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('D34').Value = '''2.47'
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('D35').Value = '''0.906'
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').Value = '1.132.81'
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('D37').Value = '''0.541'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('E39').Value = '  +0.39%  '
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('D41').Value = '''99.16'
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').Value = '''0.799'
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('E43').Value = '  +2.59%  '
$ws.Range('D44').Value = '''56.55'
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('E45').Value = '  +3.27%  '
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('E47').Value = '  +1.62%  '
$ws.Range('E48').Value = '  -0.18%  '
$ws.Range('E49').Value = '  +0.40%  '
$ws.Range('D50').Value = '''0.0952'
$ws.Range('E51').Value = '  +3.62%  '
